$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4615384615384616
$ws.Range("C2").Value = 0.6
$ws.Range("D2").Value = 0.5217391304347826

$ws.Range("B3").Value = 0.6363636363636364
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 0.5600000000000001

$ws.Range("B5").Value = 0.548951048951049
$ws.Range("C5").Value = 0.55
$ws.Range("D5").Value = 0.5408695652173914

$ws.Range("B6").Value = 0.5635198135198135
$ws.Range("D6").Value = 0.5440579710144928

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0

$ws.Range("B8").Value = 0.5833333333333334
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0.7368421052631579

$ws.Range("B9").Value = 0.5833333333333334
$ws.Range("C9").Value = 0.5833333333333334
$ws.Range("D9").Value = 0.5833333333333334
$ws.Range("E9").Value = 0.5833333333333334

$ws.Range("B10").Value = 0.2916666666666667
$ws.Range("C10").Value = 0.5
$ws.Range("D10").Value = 0.3684210526315789

$ws.Range("B11").Value = 0.3402777777777778
$ws.Range("C11").Value = 0.5833333333333334
$ws.Range("D11").Value = 0.4298245614035088

$ws.Range("B12").Value = 0.5
$ws.Range("C12").Value = 0.6
$ws.Range("D12").Value = 0.5454545454545454

$ws.Range("B13").Value = 0.6666666666666666
$ws.Range("D13").Value = 0.6153846153846153

$ws.Range("B14").Value = 0.5833333333333334
$ws.Range("C14").Value = 0.5833333333333334
$ws.Range("D14").Value = 0.5833333333333334
$ws.Range("E14").Value = 0.5833333333333334

$ws.Range("B15").Value = 0.5833333333333333
$ws.Range("C15").Value = 0.5857142857142856
$ws.Range("D15").Value = 0.5804195804195804

$ws.Range("B16").Value = 0.5972222222222222
$ws.Range("C16").Value = 0.5833333333333334
$ws.Range("D16").Value = 0.5862470862470862

$ws.Range("B17").Value = 0.4375
$ws.Range("D17").Value = 0.5384615384615384

$ws.Range("B18").Value = 0.625
$ws.Range("C18").Value = 0.3571428571428572
$ws.Range("D18").Value = 0.4545454545454545

$ws.Range("B19").Value = 0.5
$ws.Range("C19").Value = 0.5
$ws.Range("D19").Value = 0.5
$ws.Range("E19").Value = 0.5

$ws.Range("B20").Value = 0.53125
$ws.Range("C20").Value = 0.5285714285714286
$ws.Range("D20").Value = 0.4965034965034965

$ws.Range("B21").Value = 0.546875
$ws.Range("C21").Value = 0.5
$ws.Range("D21").Value = 0.4895104895104894

$ws.Range("B22").Value = 0.5384615384615384
$ws.Range("C22").Value = 0.7
$ws.Range("D22").Value = 0.608695652173913

$ws.Range("B23").Value = 0.7272727272727273
$ws.Range("D23").Value = 0.64

$ws.Range("B24").Value = 0.625
$ws.Range("C24").Value = 0.625
$ws.Range("D24").Value = 0.625
$ws.Range("E24").Value = 0.625

$ws.Range("B25").Value = 0.6328671328671329
$ws.Range("C25").Value = 0.6357142857142857
$ws.Range("D25").Value = 0.6243478260869565

$ws.Range("B26").Value = 0.6486013986013986
$ws.Range("C26").Value = 0.625
$ws.Range("D26").Value = 0.6269565217391304
